$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate the K column (column G) values, replacing the old Strike#-derived
# values with the newly computed K values.
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 2
